$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-5 (John, Jane, James, Jack) - entire row shift up
$ws.Range("A2:A5").EntireRow.Delete()

# After the shift, "Helem" row (originally row 9) is now row 5 - delete it too
$ws.Range("A5:A5").EntireRow.Delete()

$ws.Range("A7").Select()
